# Update countries & provincias Spain
# Re-sync daily COVID-19 country table: refresh case counts, re-rank a few
# countries whose totals crossed each other (so the displayed country name
# for a given row changes), and bump the 'last updated' timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Countries that swapped rank with a neighbour: update the country name
#    shown in column A for that row.
$paisRenames = @(
    @{ Row = 37; Pais = 'Oman' },
    @{ Row = 38; Pais = 'China' },
    @{ Row = 39; Pais = 'Suecia' },
    @{ Row = 124; Pais = 'Mayotte' },
    @{ Row = 125; Pais = 'Somalia' },
    @{ Row = 143; Pais = 'Aruba' },
    @{ Row = 144; Pais = 'Jamaica' },
    @{ Row = 145; Pais = 'Jordania' },
    @{ Row = 146; Pais = 'Malta' }
)
foreach ($r in $paisRenames) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Pais
}

# 2) Refreshed daily figures per row (columns B:H = Casos totales, Nuevos
#    casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes).
$datos = @(
    @{ Row = 4; CasosTotales = 6032814; NuevosCasos = 32449; CasosActivos = 3328058; Recuperados = 2520367; MuertesHoy = 736; Muertes = 184389 },
    @{ Row = 6; CasosTotales = 3384575; NuevosCasos = 76826; CasosActivos = 2583063; Recuperados = 739818; MuertesHoy = 1065; Muertes = 61694 },
    @{ Row = 8; CasosTotales = 618286; NuevosCasos = 2585; CasosActivos = 531338; Recuperados = 73320; MuertesHoy = 126; Muertes = 13628 },
    @{ Row = 12; CasosTotales = 451792 },
    @{ Row = 22; CasosTotales = 259698; NuevosCasos = 6111; CasosActivos = 85984; Recuperados = 143138; MuertesHoy = 32; Muertes = 30576 },
    @{ Row = 23; CasosTotales = 240558; NuevosCasos = 1558; Recuperados = 16966 },
    @{ Row = 32; CasosTotales = 110186; NuevosCasos = 1783; CasosActivos = 87714; Recuperados = 21589; MuertesHoy = 8; Muertes = 883 },
    @{ Row = 37; CasosTotales = 85005; NuevosCasos = 187; CasosActivos = 79608; Recuperados = 4747; MuertesHoy = 4; Muertes = 650 },
    @{ Row = 38; CasosTotales = 85004; NuevosCasos = 8; CasosActivos = 80046; Recuperados = 324; MuertesHoy = 0; Muertes = 4634 },
    @{ Row = 39; CasosTotales = 83898; NuevosCasos = 0; CasosActivos = 0; Recuperados = 0; MuertesHoy = 5; Muertes = 5820 },
    @{ Row = 64; CasosTotales = 37292; NuevosCasos = 985; CasosActivos = 14471; Recuperados = 22424; MuertesHoy = 11; Muertes = 397 },
    @{ Row = 97; CasosTotales = 8982; NuevosCasos = 46; CasosActivos = 8508; Recuperados = 416; MuertesHoy = 2; Muertes = 58 },
    @{ Row = 103; CasosTotales = 7329; NuevosCasos = 104; CasosActivos = 4691; Recuperados = 2610 },
    @{ Row = 106; CasosTotales = 6292; NuevosCasos = 41; CasosActivos = 5010; Recuperados = 1093; MuertesHoy = 10; Muertes = 189 },
    @{ Row = 108; CasosActivos = 5307; Recuperados = 16 },
    @{ Row = 117; CasosTotales = 3806; NuevosCasos = 47; CasosActivos = 3195; Recuperados = 519 },
    @{ Row = 119; CasosTotales = 3699; NuevosCasos = 69; CasosActivos = 2749; Recuperados = 912 },
    @{ Row = 124; CasosTotales = 3301; NuevosCasos = 64; CasosActivos = 2964; Recuperados = 297; MuertesHoy = 1; Muertes = 40 },
    @{ Row = 125; CasosTotales = 3275; CasosActivos = 2443; Recuperados = 737; Muertes = 95 },
    @{ Row = 127; CasosTotales = 2986; NuevosCasos = 2; CasosActivos = 2830; Recuperados = 144 },
    @{ Row = 131; CasosTotales = 2730; NuevosCasos = 13; CasosActivos = 2054; Recuperados = 550 },
    @{ Row = 133; CasosTotales = 2514; NuevosCasos = 4; Recuperados = 1177 },
    @{ Row = 134; CasosTotales = 2504; NuevosCasos = 64; CasosActivos = 569; Recuperados = 1835; MuertesHoy = 2; Muertes = 100 },
    @{ Row = 135; CasosTotales = 2415; NuevosCasos = 83; CasosActivos = 1335; Recuperados = 975; MuertesHoy = 2; Muertes = 105 },
    @{ Row = 140; CasosTotales = 2013; NuevosCasos = 10; CasosActivos = 1581; Recuperados = 362; MuertesHoy = 1; Muertes = 70 },
    @{ Row = 143; CasosTotales = 1848; NuevosCasos = 88; CasosActivos = 608; Recuperados = 1232; Muertes = 8 },
    @{ Row = 144; CasosTotales = 1804; NuevosCasos = 72; CasosActivos = 846; Recuperados = 939; Muertes = 19 },
    @{ Row = 145; CasosTotales = 1801; NuevosCasos = 45; CasosActivos = 1364; Recuperados = 422; Muertes = 15 },
    @{ Row = 146; CasosTotales = 1788; NuevosCasos = 37; CasosActivos = 1121; Recuperados = 657; Muertes = 10 },
    @{ Row = 153; CasosTotales = 1410; NuevosCasos = 38; Recuperados = 712 },
    @{ Row = 176; CasosTotales = 424; NuevosCasos = 5; Recuperados = 188 }
)
foreach ($d in $datos) {
    if ($d.ContainsKey('CasosTotales'))  { $ws.Cells.Item($d.Row, 2).Value = $d.CasosTotales }
    if ($d.ContainsKey('NuevosCasos'))   { $ws.Cells.Item($d.Row, 3).Value = $d.NuevosCasos }
    if ($d.ContainsKey('CasosActivos'))  { $ws.Cells.Item($d.Row, 4).Value = $d.CasosActivos }
    if ($d.ContainsKey('Recuperados'))   { $ws.Cells.Item($d.Row, 5).Value = $d.Recuperados }
    if ($d.ContainsKey('CasosCriticos')) { $ws.Cells.Item($d.Row, 6).Value = $d.CasosCriticos }
    if ($d.ContainsKey('MuertesHoy'))    { $ws.Cells.Item($d.Row, 7).Value = $d.MuertesHoy }
    if ($d.ContainsKey('Muertes'))       { $ws.Cells.Item($d.Row, 8).Value = $d.Muertes }
}

# 3) Bump the 'last updated' footer timestamp.
$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 27 de Agosto de 2020 a las 22:19'

